# "adicionado politica de preco"
#
# The sheet gains two new columns inserted after "preco" (B) and before the
# old "full" column (old C): "modelo" (C) and "politica" (D).
# The previously existing columns full/tipo/link shift two columns to the
# right (old C/D/E -> new E/F/G). For every data row the "tipo" values are
# also normalised to lower-case, and the link's tracking_id query value is
# updated to a new GUID.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) -----------------------------------------------
# Grab the existing bold/centered header style from A1 so the two brand
# new header cells (C1, D1) and the shifted ones (F1, G1) keep looking
# like headers.
$ws.Range("A1").Copy()
$ws.Range("C1:G1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C1").Value = "modelo"
$ws.Range("D1").Value = "politica"
$ws.Range("E1").Value = "full"
$ws.Range("F1").Value = "tipo"
$ws.Range("G1").Value = "link"

# ---- Data rows (rows 2-25) ---------------------------------------------
$oldTrackingId = "6fcea338-e789-4443-9675-16c3b4a01a1a"
$newTrackingId = "c1f14ab8-1380-4955-b0e4-59204e5b1141"

for ($row = 2; $row -le 25; $row++) {
    # Capture the old full/tipo/link values before any column gets
    # overwritten.
    $oldFull = $ws.Cells.Item($row, 3).Value2   # old C: "full"
    $oldTipo = $ws.Cells.Item($row, 4).Value2   # old D: "tipo"
    $oldLink = $ws.Cells.Item($row, 5).Value2   # old E: "link"

    $newLink = $oldLink.Replace($oldTrackingId, $newTrackingId)
    $newTipo = $oldTipo.ToLower()

    # Write the shifted columns first (G, F, E) while the source columns
    # (C, D, E) still hold their original values.
    $ws.Cells.Item($row, 7).Value = $newLink     # G: link (new tracking id)
    $ws.Cells.Item($row, 6).Value = $newTipo     # F: tipo (lower-cased)
    $ws.Cells.Item($row, 5).Value = $oldFull     # E: full (unchanged, "NA")

    # Now populate the two brand-new columns.
    $ws.Cells.Item($row, 4).Value = ""           # D: politica (blank)
    $ws.Cells.Item($row, 3).Value = "Sem Modelo" # C: modelo
}

Write-Output "politica de preco columns added"
